$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductsTests")

$ws.Range("A4").Value = "TC_011"
$ws.Range("B4").Value = "7.99, 9.99, 15.99, 15.99, 29.99, 49.99"

$ws.Range("A5").Value = "TC_012"
$ws.Range("B5").Value = "49.99, 29.99, 15.99, 15.99, 9.99, 7.99"

$ws.Activate()
$ws.Range("B6").Select()
